# edit.ps1 -- apply the "AI essay" -> "Arts in society" essay rewrite
# described by the reference diff, via Word COM-interop (Find/Replace +
# a couple of InsertAfter calls for brand-new sentences and a trailing
# empty paragraph).

$d = $word.ActiveDocument

function Replace-Text($oldText, $newText) {
    $r = $d.Content
    $ok = $r.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $ok) {
        Write-Output "MISSING: $oldText"
    }
    return $ok
}

function Append-After($oldText, $extraText) {
    # Find oldText (post-replacement), collapse to its end, and insert
    # extraText right after it -- used where the diff adds brand new
    # sentences following an existing (possibly just-replaced) run.
    $r = $d.Content
    $ok = $r.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $ok) {
        Write-Output "MISSING-APPEND-ANCHOR: $oldText"
        return
    }
    $r.Collapse(0)
    $r.InsertAfter($extraText)
}

# ---------------------------------------------------------------------
# Title
# ---------------------------------------------------------------------
Replace-Text "The Allure of Artificial Intelligence" "Understanding the Role of the Arts in Society"

# ---------------------------------------------------------------------
# Author name: "Sarah Jones" -> "Amelia J" + "." + " Clayton"
# ---------------------------------------------------------------------
Replace-Text "Sarah Jones" "Amelia J"
Append-After "Amelia J" "."
Append-After "Amelia J." " Clayton"

# ---------------------------------------------------------------------
# Email
# ---------------------------------------------------------------------
Replace-Text "username@xyzdomain" "claytonamelia123@gmail"

# ---------------------------------------------------------------------
# Body paragraph 1 (intro)
# ---------------------------------------------------------------------
Replace-Text "Artificial intelligence (AI) has become a captivating field that continues to fascinate and challenge our understanding of technology" "1"

Replace-Text " Its influence is growing in diverse domains, reshaping industries and transforming our daily lives" " Exploring the world of arts allows individuals to step into the realms of imagination, creativity, and self-expression"

Replace-Text " This essay aims to explore the allure of AI and dive into the key aspects that make it such a captivating field of study" " It's a multifaceted realm where various art forms, from visual masterpieces to enchanting melodies, contribute to a broader narrative of our shared existence"

Replace-Text " We will embark on a journey to comprehend its impact on various sectors, examine the underlying principles powering AI's capabilities, and contemplate the ethical considerations surrounding its development and deployment" " Whether it's the vibrant canvas of a painting, the profound lyrics of a song, or the eloquence of a narrative, art enables us to transcend boundaries, engage with our emotions, and embark on a transformative journey of self-discovery"
Append-After " Whether it's the vibrant canvas of a painting, the profound lyrics of a song, or the eloquence of a narrative, art enables us to transcend boundaries, engage with our emotions, and embark on a transformative journey of self-discovery" "."
Append-After " Whether it's the vibrant canvas of a painting, the profound lyrics of a song, or the eloquence of a narrative, art enables us to transcend boundaries, engage with our emotions, and embark on a transformative journey of self-discovery." " In this exploration, we'll delve into the essence of the arts and its multifaceted role in shaping human experiences and societies"

# ---------------------------------------------------------------------
# Body paragraph 1, second block (after first <br/><br/>)
# ---------------------------------------------------------------------
Replace-Text "Humans have always strived to understand and replicate intelligence" "2"

Replace-Text " The pursuit of creating intelligent machines dates back centuries, with AI emerging as the frontier of this quest" " Art isn't merely about aesthetics; it embodies a tapestry of cultural narratives, historical contexts, and societal values"

Replace-Text " Its ability to learn, reason, and make decisions autonomously has captivated researchers, scientists, and engineers worldwide" " Through studying and understanding the arts, we delve into the depths of human existence"

Replace-Text " AI's potential to solve complex problems and automate tasks considered challenging or impossible for humans has drawn immense interest and enthusiasm" " Whether it's the intricate design of ancient artifacts speaking of forgotten civilizations or the profound narratives of literary classics shedding light on complex human conditions, art acts as a medium to connect across time and cultures"

Replace-Text " From healthcare and finance to autonomous vehicles and space exploration, AI's applications are far-reaching and continue to expand" " In engaging with artistic creations, we can develop empathy, fostering a greater appreciation for diverse perspectives and a more inclusive world"

# ---------------------------------------------------------------------
# Body paragraph 1, third block (after second <br/><br/>)
# ---------------------------------------------------------------------
Replace-Text "AI's allure lies not only in its practical applications but also in its profound implications for our understanding of intelligence itself" "3"

Replace-Text " The study of AI has shed light on the intricate mechanisms of human cognition and behavior, challenging traditional notions of consciousness and self-awareness" " Beyond its intrinsic value, the arts have far-reaching societal impacts, contributing to economic vibrancy, social harmony, and individual well-being"

Replace-Text " As AI systems become more sophisticated, we are confronted with questions about the nature of intelligence, the limits of computation, and the potential for machines to surpass human capabilities" " The arts industry employs millions worldwide, fueling economic growth and providing employment opportunities"

Replace-Text " These intellectual explorations drive the pursuit of AI forward, fueling a vibrant and rapidly evolving field of research" " By promoting cultural tourism and fostering creativity, the arts attract visitors and enhance community vitality"
Append-After " By promoting cultural tourism and fostering creativity, the arts attract visitors and enhance community vitality" "."
Append-After " By promoting cultural tourism and fostering creativity, the arts attract visitors and enhance community vitality." " Moreover, the arts can serve as a powerful tool for education, promoting critical thinking, enhancing communication skills, and cultivating creativity"
Append-After " Moreover, the arts can serve as a powerful tool for education, promoting critical thinking, enhancing communication skills, and cultivating creativity" "."
Append-After " Moreover, the arts can serve as a powerful tool for education, promoting critical thinking, enhancing communication skills, and cultivating creativity." " When individuals participate in artistic activities, they develop crucial life skills that extend beyond academic settings, equipping them for success in various endeavors"

# ---------------------------------------------------------------------
# Summary paragraph
# ---------------------------------------------------------------------
Replace-Text "The allure of Artificial Intelligence stems from its far-reaching impact, its profound implications for our understanding of intelligence itself, and its transformative potential across diverse fields" "In this essay, we explored the multifaceted role of the arts in society"

Replace-Text " AI's ability to learn, reason, and make decisions autonomously has opened up new possibilities for solving complex problems and automating tasks" " We discussed how the arts foster creativity, enable self-expression, bridge cultural gaps, and nurture empathy"

Replace-Text " Its influence is shaping industries, transforming our daily lives, and challenging our notions of intelligence and consciousness" " Moreover, we highlighted the instrumental role the arts play in education, economic development, and communal well-being"

Replace-Text " While the development and deployment of AI raise ethical and societal considerations, the allure of this captivating field continues to drive innovation and exploration at the forefront of technology" " By appreciating and engaging with the arts in all its forms, we not only enrich our individual lives but also contribute to the vibrancy and progress of our communities"

# ---------------------------------------------------------------------
# New trailing empty paragraph after the Summary paragraph
# ---------------------------------------------------------------------
$d.Paragraphs.Add() | Out-Null

Write-Output "done"
